$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(new C value, new D value)
$updates = @{
    2 = @(41589, 60048447)
    3 = @(98772, 144683350)
    4 = @(33624, 49768313)
    5 = @(9616, 14284242)
    6 = @(2379, 3534471)
    7 = @(249, 368593)
    12 = @(44732, 60515782)
    13 = @(10534, 15212918)
    14 = @(27819, 40765897)
    15 = @(8818, 13085166)
    16 = @(2335, 3467875)
    20 = @(10976, 14449504)
    21 = @(14479, 20865750)
    22 = @(33662, 49355296)
    23 = @(10835, 16099239)
    24 = @(2876, 4274115)
    26 = @(44, 65453)
    27 = @(12503, 16615544)
    28 = @(8423, 12181684)
    29 = @(24221, 35524253)
    30 = @(8334, 12389707)
    31 = @(2127, 3172708)
    32 = @(430, 635621)
    34 = @(8992, 11844848)
    35 = @(3684, 5314068)
    36 = @(8631, 12611045)
    37 = @(3404, 5047008)
    38 = @(875, 1303555)
    39 = @(181, 269186)
    41 = @(2761, 3715647)
    42 = @(18834, 27191903)
    43 = @(54807, 80290965)
    44 = @(20098, 29832081)
    46 = @(1411, 2106144)
    50 = @(18270, 24186994)
    51 = @(2370, 3438526)
    52 = @(7981, 11723028)
    53 = @(2675, 3992633)
    54 = @(843, 1259414)
    57 = @(7953, 10941837)
    58 = @(1650, 3316440)
    59 = @(3939, 7882312)
    60 = @(1557, 3121294)
    61 = @(528, 1053083)
    62 = @(191, 394487)
    64 = @(2561, 4752864)
    65 = @(16976, 24498817)
    66 = @(48309, 70617001)
    67 = @(16882, 25083962)
    69 = @(1091, 1622199)
    73 = @(16130, 21164772)
    74 = @(62563, 90931790)
    75 = @(170320, 250643438)
    76 = @(72772, 108376532)
    77 = @(23812, 35562074)
    78 = @(6195, 9244877)
    79 = @(451, 671458)
    85 = @(61562, 83003511)
    86 = @(5133, 7435873)
    87 = @(12616, 18525814)
    88 = @(4136, 6163126)
    89 = @(1455, 2173611)
    93 = @(5874, 7877874)
    94 = @(1849, 2661993)
    95 = @(5865, 8642757)
    96 = @(2109, 3143431)
    98 = @(223, 339109)
    99 = @(25, 37500)
    101 = @(3984, 5280363)
    102 = @(919, 1779061)
    103 = @(615, 1259284)
    107 = @(11900, 17245765)
    108 = @(31121, 45675276)
    109 = @(10429, 15502760)
    110 = @(2891, 4309571)
    115 = @(10493, 13803149)
    116 = @(33414, 48133995)
    117 = @(70927, 103732868)
    118 = @(22744, 33782586)
    119 = @(6526, 9715052)
    120 = @(1305, 1950233)
    121 = @(127, 185295)
    125 = @(27695, 36870555)
    126 = @(39919, 57542589)
    127 = @(83348, 121766362)
    128 = @(25527, 37877249)
    129 = @(6967, 10353009)
    130 = @(1476, 2187261)
    134 = @(34211, 45291943)
    135 = @(14518, 20998973)
    136 = @(34618, 50805744)
    137 = @(12211, 18142771)
    138 = @(3232, 4817875)
    143 = @(11603, 15414361)
    144 = @(39222, 56622636)
    145 = @(89798, 131454155)
    146 = @(26835, 39863611)
    147 = @(7154, 10660774)
    148 = @(1750, 2600406)
    151 = @(31717, 42621580)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}
